$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.044.00"
$ws.Range("E2").Value = "  +1.92%  "

$ws.Range("D3").Value = "3.895.65"
$ws.Range("E3").Value = "  +0.47%  "

$ws.Range("D5").Value = "'484.47"
$ws.Range("E5").Value = "  +2.77%  "

$ws.Range("D6").Value = "'145.20"
$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").Value = "'0.619"
$ws.Range("E7").Value = "  -1.53%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("D9").Value = "'0.721"
$ws.Range("E9").Value = "  -3.33%  "

$ws.Range("D10").Value = "'0.166"
$ws.Range("E10").Value = "  +4.59%  "

$ws.Range("D11").Value = "'0.0000360"
$ws.Range("E11").Value = "  +12.83%  "

$ws.Range("D12").Value = "'42.65"
$ws.Range("E12").Value = "  -2.02%  "

$ws.Range("D13").Value = "'10.61"
$ws.Range("E13").Value = "  +1.06%  "

$ws.Range("D14").Value = "4.506.12"
$ws.Range("E14").Value = "  +0.23%  "

$ws.Range("D15").Value = "'14.63"
$ws.Range("E15").Value = "  -2.45%  "

$ws.Range("D16").Value = "3.904.17"
$ws.Range("E16").Value = "  +1.29%  "

$ws.Range("E17").Value = "  -0.36%  "

$ws.Range("D18").Value = "'19.72"
$ws.Range("E18").Value = "  -2.55%  "

$ws.Range("E19").Value = "  -3.90%  "

$ws.Range("D20").Value = "68.138.60"
$ws.Range("E20").Value = "  +1.53%  "

$ws.Range("D21").Value = "'435.77"
$ws.Range("E21").Value = "  +0.98%  "

$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D22").Value = "'14.62"
$ws.Range("E22").Value = "  -2.78%  "

$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D23").Value = "'3.33"
$ws.Range("E23").Value = "  -0.83%  "

$ws.Range("D24").Value = "'88.55"
$ws.Range("E24").Value = "  +0.59%  "

$ws.Range("D25").Value = "'11.64"
$ws.Range("E25").Value = "  +15.09%  "

$ws.Range("E26").Value = "  +0.51%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'38.41"
$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'10.33"
$ws.Range("E28").Value = "  +2.76%  "

$ws.Range("D29").Value = "'5.79"
$ws.Range("E29").Value = "  +3.96%  "

$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "'13.42"
$ws.Range("E30").Value = "  -3.68%  "

$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").Value = "'690.60"
$ws.Range("E31").Value = "  -5.79%  "

$ws.Range("E32").Value = "  -2.56%  "

$ws.Range("D33").Value = "'2.85"
$ws.Range("E33").Value = "  +2.40%  "

$ws.Range("D34").Value = "0.0₃0932"
$ws.Range("E34").Value = "  +29.16%  "

$ws.Range("D35").Value = "'41.24"
$ws.Range("E35").Value = "  -5.52%  "

$ws.Range("D36").Value = "'59.09"
$ws.Range("E36").Value = "  +0.31%  "

$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.150"
$ws.Range("E37").Value = "  -7.15%  "

$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'5.65"
$ws.Range("E38").Value = "  +3.23%  "

$ws.Range("E39").Value = "  -0.08%  "

$ws.Range("D40").Value = "'0.0471"
$ws.Range("E40").Value = "  -3.00%  "

$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").Value = "'2.71"
$ws.Range("E41").Value = "  +3.25%  "

$ws.Range("B42").Value = "ThetaToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D42").Value = "'2.98"
$ws.Range("E42").Value = "  -1.08%  "

$ws.Range("D43").Value = "'2.96"
$ws.Range("E43").Value = "  +6.70%  "

$ws.Range("D44").Value = "'0.346"
$ws.Range("E44").Value = "  +0.14%  "

$ws.Range("D45").Value = "'0.141"
$ws.Range("E45").Value = "  -1.17%  "

$ws.Range("E46").Value = "  -0.07%  "

$ws.Range("E47").Value = "  -2.19%  "

$ws.Range("E48").Value = "  -3.22%  "

$ws.Range("D49").Value = "'146.22"
$ws.Range("E49").Value = "  +2.45%  "

$ws.Range("E50").Value = "  -2.24%  "

$ws.Range("D51").Value = "'2.82"
$ws.Range("E51").Value = "  -3.69%  "

